$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 47. This shifts the existing rows 47-58
# down to 48-59 (content, styles, everything moves with them), exactly
# matching the "row inserted above" pattern seen in the target diff.
$ws.Rows("47:47").Insert()

# The newly inserted row inherits the formatting of the row above it
# (row 46, style index 4). The target row instead needs the "s=3"
# styling used by neighbouring GMHO rows, so copy that formatting from
# an existing s="3" row (row 49, after the shift) onto the whole new
# row 47.
$fmtSource = $ws.Range("A49:V49")
$fmtDest = $ws.Range("A47:V47")
$fmtSource.Copy()
$fmtDest.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row's content ("psychological treatment" / GMHO:0000239).
$ws.Range("A47").Value = "GMHO:0000239"
$ws.Range("B47").Value = "psychological treatment"
$ws.Range("C47").Value = "Mental health intervention content that uses communication or recommended tasks to assess and improve a person’s adaptive mental or behavioural functioning."
$ws.Range("D47").Value = "mental health intervention content"
$ws.Range("K47").Value = "https://bciosearch.org/BCIO_050364"
$ws.Range("L47").Value = "BCIO:050364"
$ws.Range("P47").Value = "LSR 2"
$ws.Range("S47").Value = "Proposed"
$ws.Range("V47").Value = "MS"
